$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks so we can recreate them pointing at the
# (row-shifted) cells once the new row has been inserted.
$ws.Hyperlinks.Delete()

# Insert a new row at 11 ("Liver" / "Perihepatic reactive lymph nodes" / ...),
# pushing every row from 11 downward down by one.
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = "Liver"
$ws.Range("B11").Value = "Perihepatic reactive lymph nodes"
$ws.Range("C11").Value = "Clip 1-Bmode"
$ws.Range("D11").Value = "https://youtu.be/kaROVVBl9Bc"
$ws.Range("D11").Style = "Collegamento ipertestuale"

# Re-create the hyperlinks at their new (shifted where applicable) locations.
$ws.Hyperlinks.Add($ws.Range("D3"), "https://youtu.be/zxTC0YBY2RY")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://youtu.be/xBfd04F4Ni8")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://youtu.be/91M82AIMyu0")
$ws.Hyperlinks.Add($ws.Range("D18"), "https://youtu.be/qushjTAy6XQ")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://youtu.be/pc-vbxSRTbs")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://youtu.be/DjI1kEnzfSQ")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://youtu.be/JvwODCASLYQ")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://youtu.be/U3ydTsRwxok")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://youtu.be/15o_Km86IzM")

# Re-apply the hyperlink style to cells whose style may have been disturbed.
$ws.Range("D3").Style = "Collegamento ipertestuale"
$ws.Range("D8").Style = "Collegamento ipertestuale"
$ws.Range("D10").Style = "Collegamento ipertestuale"
$ws.Range("D12").Style = "Collegamento ipertestuale"
$ws.Range("D13").Style = "Collegamento ipertestuale"
$ws.Range("D15").Style = "Collegamento ipertestuale"
$ws.Range("D16").Style = "Collegamento ipertestuale"
$ws.Range("D17").Style = "Collegamento ipertestuale"
$ws.Range("D18").Style = "Collegamento ipertestuale"

# Sort stays the same, but refresh the view selection/scroll like Excel would
# after editing cell D11.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D11").Select()
